$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 660
$ws.Range("F5").Value = 641
$ws.Range("F6").Value = 1434
$ws.Range("F10").Value = 2091
$ws.Range("F15").Value = 234
$ws.Range("F16").Value = 1035
$ws.Range("F17").Value = 404
$ws.Range("F18").Value = 64
$ws.Range("F19").Value = 153
$ws.Range("F20").Value = 4027
$ws.Range("F21").Value = 1222
$ws.Range("F22").Value = 3168
$ws.Range("F23").Value = 311
$ws.Range("F24").Value = 97
$ws.Range("F25").Value = 2937
$ws.Range("F26").Value = 2937
$ws.Range("F27").Value = 4595
$ws.Range("F29").Value = 952
$ws.Range("F31").Value = 3018
$ws.Range("F32").Value = 288
$ws.Range("F34").Value = 108
$ws.Range("F39").Value = 96
$ws.Range("F40").Value = 1213
$ws.Range("F41").Value = 785
$ws.Range("F45").Value = 40
$ws.Range("F46").Value = 196
$ws.Range("F50").Value = 3663

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 958
$ws.Range("F21").Value = 37
$ws.Range("F23").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1593

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 641
$ws.Range("F5").Value = 1434
$ws.Range("F10").Value = 2091
$ws.Range("F14").Value = 958
$ws.Range("F17").Value = 234
$ws.Range("F18").Value = 1035
$ws.Range("F20").Value = 404
$ws.Range("F21").Value = 153
$ws.Range("F22").Value = 4027
$ws.Range("F24").Value = 1222
$ws.Range("F26").Value = 3168
$ws.Range("F27").Value = 2937
$ws.Range("F28").Value = 2937
$ws.Range("F29").Value = 4595
$ws.Range("F30").Value = 952
$ws.Range("F31").Value = 3018
$ws.Range("F32").Value = 288
$ws.Range("F33").Value = 108
$ws.Range("F37").Value = 96
$ws.Range("F38").Value = 1213
$ws.Range("F39").Value = 785
$ws.Range("F43").Value = 37
$ws.Range("F44").Value = 40
$ws.Range("F46").Value = 196
$ws.Range("F50").Value = 3663
